$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("选文信息")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.theatlantic.com/business/archive/2016/06/would-a-world-without-work-be-so-bad/488711/")
